$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 110; existing rows 110..178 shift down to 111..179.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44518
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112043
$ws.Cells.Item(110, 7).Value = "Pepino ensalada"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 200
$ws.Cells.Item(110, 11).Value = 11000
$ws.Cells.Item(110, 12).Value = 11000
$ws.Cells.Item(110, 13).Value = 11000
$ws.Cells.Item(110, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(110, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 16).Value = 183
$ws.Cells.Item(110, 17).Value = 60
$ws.Cells.Item(110, 18).Value = "Hortaliza"
